$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 71; this shifts existing rows 71-183 down to 72-184.
$ws.Rows(71).Insert()

# Populate the new row 71 with a new weekly record (same as the old row 71
# record except for the date and volume, which are updated).
$ws.Range("A71").Value = 4
$ws.Range("B71").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C71").Value = "Los Lagos"
$ws.Range("D71").Value = 44533
$ws.Range("E71").Value = 10
$ws.Range("F71").Value = "Fruta"
$ws.Range("G71").Value = 100102
$ws.Range("H71").Value = "Cítricos"
$ws.Range("I71").Value = 100102006
$ws.Range("J71").Value = "Pomelo"
$ws.Range("K71").Value = "Start Ruby"
$ws.Range("L71").Value = "Primera"
$ws.Range("M71").Value = 200
$ws.Range("N71").Value = 11000
$ws.Range("O71").Value = 12000
$ws.Range("P71").Value = 11500
$ws.Range("Q71").Value = "$/caja 14 kilos empedrada"
$ws.Range("R71").Value = "Región de O'Higgins"
$ws.Range("S71").Value = 821
$ws.Range("T71").Value = 14

$ws.Range("D71").NumberFormat = "YYYY-MM-DD HH:MM:SS"
